# Apply targeted corrections to column F (dSF) values as part of a
# "repull data, push all data, mean calculation" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$updates = @{
    "F6"  = 4
    "F10" = -5
    "F12" = -2
    "F13" = 4
    "F17" = -6
    "F25" = -7
    "F27" = -1
    "F31" = 0
    "F34" = -2
    "F41" = 2
    "F45" = 0
    "F51" = -1
    "F61" = 2
    "F66" = 2
    "F69" = -2
    "F73" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
